# CHANGED: Card Lists - Added starship class abbreviations and changed Purple Wing starship names.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Blue Wing Battleship flagships -> add "TRB " prefix
$ws.Range("B14").Value = "TRB Ardor"
$ws.Range("B15").Value = "TRB Faith"
$ws.Range("B16").Value = "TRB Justice"
$ws.Range("B17").Value = "TRB Truth"

# Dra'tar starships -> add class-abbreviation prefixes
$ws.Range("B25").Value = "DB Agony"
$ws.Range("B28").Value = "DSF Bane"
$ws.Range("B30").Value = "DAF Cataclysm"
$ws.Range("B31").Value = "DSF Curse"
$ws.Range("B32").Value = "DSF Decay"
$ws.Range("B33").Value = "DAC Disaster"
$ws.Range("B34").Value = "DB Horror"

# Green Wing Dreadnoughts -> renamed class to Frigate, names get "TRF " prefix
$ws.Range("B48").Value = "TRF Moscow"
$ws.Range("J48").Value = "Frigate"
$ws.Range("B49").Value = "TRF Washington"
$ws.Range("J49").Value = "Frigate"

# Purple Wing Interdictors -> renamed (new names), "TRI " prefix
$ws.Range("B61").Value = "TRI Pegasus"
$ws.Range("B62").Value = "TRI Golem"
